$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Seed the 16 brand-new shared strings in the exact order the original
#    authoring session created them (si index 74-89), so the shared string
#    table lines up with the target workbook.
$ws.Range("J164").Value = "* Radu"
$ws.Range("J167").Value = "* Tolga"
$ws.Range("J168").Value = "Design x axis"
$ws.Range("J169").Value = "encoder holder for x axis"
$ws.Range("J171").Value = "* Nikkita"
$ws.Range("J172").Value = "Test y axis"
$ws.Range("J174").Value = "* Tobias "
$ws.Range("J175").Value = "Battery saftey"
$ws.Range("J146").Value = "Nothing yet"
$ws.Range("J177").Value = "* Mandvias"
$ws.Range("J178").Value = "Code for equation in C"
$ws.Range("J165").Value = "Battery Bus (Electrical Suply Hub)"
$ws.Range("J151").Value = "Translated formula into a function is C code (during meeting)"
$ws.Range("J155").Value = "Designing a holder for the motor (during meeting)"
$ws.Range("J159").Value = "nothin'"
$ws.Range("J161").Value = "·        Electrical Pump"

# 2) Fill in the remaining new cells that reuse pre-existing shared strings
#    or hold literal numbers -- order doesn't matter for these.
$ws.Range("J144").Value = 45789
$ws.Range("J145").Value = "·        Pump"
$ws.Range("J149").Value = " "
$ws.Range("J150").Value = "·        Screen"
$ws.Range("J154").Value = "·        X axis motor"
$ws.Range("J158").Value = "·        Y axis motor"
$ws.Range("J162").Value = "nothin'"

# 3) Apply styles by copying formats from existing reference cells
#    (J18=date style s=11, J19=symbol-bullet style s=10, J20=Aptos style s=9).
#    This reuses the existing style table instead of synthesising new xf ids.
#    NOTE: multi-area (comma-joined) ranges silently drop every area but the
#    first when used with PasteSpecial, so each destination is pasted alone.
foreach ($addr in @("J146", "J147", "J148", "J149", "J151", "J152", "J153", "J155", "J156", "J157", "J159", "J160", "J161")) {
  $ws.Range("J20").Copy() | Out-Null
  $ws.Range($addr).PasteSpecial(-4122)
}

foreach ($addr in @("J145", "J150", "J154", "J158")) {
  $ws.Range("J19").Copy() | Out-Null
  $ws.Range($addr).PasteSpecial(-4122)
}

foreach ($addr in @("J144")) {
  $ws.Range("J18").Copy() | Out-Null
  $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# 4) New font (Shruti) used only by the large banner cell J162.
$ws.Range("J162").Font.Name = "Shruti"
$ws.Range("J162").Font.Size = 11
$ws.Range("J162").EntireRow.RowHeight = 18.6

# 5) Selection / scroll position to match the saved view.
$excel.ActiveWindow.ScrollRow = 153
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J179").Select()
